# CVDLS-209: Fix Viral Result test to use Tube IDs instead of generated Specimen IDs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "TubeQPCRResults0001"
$ws.Range("A3").Value = "TubeQPCRResults0002"
$ws.Range("A4").Value = "TubeQPCRResults0003"
$ws.Range("A5").Value = "TubeQPCRResults0004"

# Reset the view: scroll back to column A (was topLeftCell="D1") and move
# the active selection from J6 to A6.
$ws.Range("A6").Select()
